# Landscaping Data.xlsx update
#  1. Correct the "Low" reading (column D) for 2025-07-25 (rows 534-540)
#     from 73 to 72 (Temp_Diff in column F is a shared formula and will
#     recalculate automatically).
#  2. Append 7 new observations for 2025-07-26 (rows 541-547).
#  3. Update the worksheet view (scroll position / selection) to match
#     where the user ended up after entering the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix D534:D540 (73 -> 72) -------------------------------------
for ($r = 534; $r -le 540; $r++) {
    $ws.Cells.Item($r, 4).Value = 72
}

# --- 2. Add new rows 541-547 ------------------------------------------
# Give the new date cells (column A) the same date number-format style
# already used by the rest of column A, by copying the format from the
# last existing date cell before typing the new values in.
$ws.Cells.Item(540, 1).Copy()
$ws.Range("A541:A547").PasteSpecial(-4122)

# row, Plant_Type, Plant_Size, Low, High, Rain, Growth, Pruned, Quadrant, Shade, UV, Humidity, Dew_Point, Pressure, Wind_Gust, Cloud_Cover, Visibility, AQI, Pollen
$newRows = @(
    @(541, "Flowering",    "Large",  73, 89, 0.31, 0.3,                 "No", 2, "Bright",  5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26),
    @(542, "Nonflowering", "Medium", 73, 89, 0.31, 0.3,                 "No", 3, "Neutral", 5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26),
    @(543, "Nonflowering", "Small",  73, 89, 0.31, 0.35,                "No", 3, "Neutral", 5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26),
    @(544, "Nonflowering", "Medium", 73, 89, 0.31, 0.4,                 "No", 3, "Bright",  5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26),
    @(545, "Nonflowering", "Medium", 73, 89, 0.31, 0.2,                 "No", 3, "Bright",  5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26),
    @(546, "Nonflowering", "Large",  73, 89, 0.31, 0.55000000000000004, "No", 4, "Neutral", 5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26),
    @(547, "Tree",         "Medium", 73, 89, 0.31, 1.7,                 "No", 1, "Neutral", 5, 0.79, 74, 30.11, 10, 0.37, 8.1, 51, 26)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 45864          # A: Date (7/26/2025)
    $ws.Cells.Item($r, 2).Value = $row[1]        # B: Plant_Type
    $ws.Cells.Item($r, 3).Value = $row[2]        # C: Plant_Size
    $ws.Cells.Item($r, 4).Value = $row[3]        # D: Low
    $ws.Cells.Item($r, 5).Value = $row[4]        # E: High
    # F: Temp_Diff -- filled in below as a shared formula
    $ws.Cells.Item($r, 7).Value = $row[5]        # G: Rain
    $ws.Cells.Item($r, 8).Value = $row[6]        # H: Growth
    $ws.Cells.Item($r, 9).Value = $row[7]        # I: Pruned
    $ws.Cells.Item($r, 10).Value = $row[8]       # J: Quadrant
    $ws.Cells.Item($r, 11).Value = $row[9]       # K: Shade
    $ws.Cells.Item($r, 12).Value = $row[10]      # L: UV
    $ws.Cells.Item($r, 13).Value = $row[11]      # M: Humidity
    $ws.Cells.Item($r, 14).Value = $row[12]      # N: Dew_Point
    $ws.Cells.Item($r, 15).Value = $row[13]      # O: Pressure
    $ws.Cells.Item($r, 16).Value = $row[14]      # P: Wind_Gust
    $ws.Cells.Item($r, 17).Value = $row[15]      # Q: Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $row[16]      # R: Visibility
    $ws.Cells.Item($r, 19).Value = $row[17]      # S: AQI
    $ws.Cells.Item($r, 20).Value = $row[18]      # T: Pollen
}

# F541:F547 -- ABS(Low-High), entered as one block so it is stored as a
# single shared formula, same as the rest of column F.
$ws.Range("F541:F547").Formula = "=ABS(D541-E541)"

# --- 3. Update the view: scroll down and select the new Pollen values -
$excel.ActiveWindow.ScrollRow = 523
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("T541:T547").Select()
